# Apply the parameter-file update described in the commit:
# "add parameter comparison code, change base and v4 parameterizations"
#
# Data changes on Sheet1:
#  - Friant-Kern Canal (row 13)            : C13,D13  2 -> 1
#  - Irrigation/water districts (row 14)   : C14,D14  2 -> 1
#  - Water Rights Division (SWRCB) (row 16): C16,D16  1 -> -1
#                                             E16,F16  1 -> (cleared)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1

$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1

$ws.Range("C16").Value = -1
$ws.Range("D16").Value = -1
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

# Cosmetic: move the view / selection the way the author last left the sheet.
$ws.Range("F16").Select()
$excel.ActiveWindow.ScrollRow = 12
